# Updated cryptos list on Fri Mar  1 09:26:43 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.198.23"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.65%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.426.65"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -1.04%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.08%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'408.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.73%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'133.06"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +1.89%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  -0.02%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -0.08%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.677"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -1.93%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.123"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'42.21"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -3.73%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  -1.51%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'3.967.36"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -1.16%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'20.00"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -1.56%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'8.46"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -3.30%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'3.427.97"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -1.23%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'62.137.58"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.70%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'  -2.73%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'11.03"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.27%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.0000132"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = "'3.21"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -4.54%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'85.17"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +4.04%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'315.50"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.68%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'12.79"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -3.39%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  -2.96%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  +9.58%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -2.61%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'8.26"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +1.69%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'7.70"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -1.16%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'2.73"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +3.58%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -2.18%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -4.47%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'42.73"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -4.90%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  -4.33%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  -0.05%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.0486"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -2.26%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'52.18"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -0.94%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.998"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +0.07%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  -3.85%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'2.99"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -0.94%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'2.01"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -0.23%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'137.95"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +0.40%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  -0.97%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  +0.35%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'3.99"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -0.31%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'16.90"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -4.82%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  -3.23%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'21.53"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -4.62%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'2.131.89"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -5.33%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  -4.14%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  +0.39%  "
$ws.Range("E51").Style = "Normal"
